$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Description (B2) and User Story (B4) values which were
# previously mismatched.
$desc = $ws.Range("B2").Value2
$story = $ws.Range("B4").Value2
$ws.Range("B2").Value = $story
$ws.Range("B4").Value = $desc

# Update the slot-card label text to reflect the new max limit (6 -> 18).
$ws.Range("N10").Value = "Other Slot Cards  (1 of 18)"
$ws.Range("N11").Value = "Other Slot Cards  (1 of 18)"

# Update the active selection/view state.
$ws.Range("J10").Select()
